$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D2:D7 switch from text labels ("2 (PD0)", "3 (PD1)", ...) to
# plain numeric pin numbers.
$ws.Range("D2").Value = 2
$ws.Range("D3").Value = 3
$ws.Range("D4").Value = 4
$ws.Range("D5").Value = 5
$ws.Range("D6").Value = 6
$ws.Range("D7").Value = 12

# Selection moves from E7 to D7.
$ws.Range("D7").Select()
